# Sources-And-Credits.xlsx update
# Adds 8 new sound-credit rows (freesound.org sources) before the trailing
# "Paid" / font-credit block, shifting that block down from rows 16-18 to
# rows 24-26, and re-establishes all hyperlinks afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Make room: insert 8 blank rows at row 15 (pushes old rows 15-18 -> 23-26)
# ---------------------------------------------------------------------
for ($i = 0; $i -lt 8; $i++) {
    $ws.Rows.Item(15).Insert()
}

# Inserting rows carries down formatting from the row above (row 14) onto
# the new blank rows; wipe that so the new rows start from a clean slate.
$ws.Range("A15:D22").Clear()

# ---------------------------------------------------------------------
# 2. Fill in the new credit rows
# ---------------------------------------------------------------------
$newRows = @(
    @{ Row = 15; Name = "60013__qubodup__whoosh";                                  Url = "https://freesound.org/people/qubodup/sounds/60013/" },
    @{ Row = 16; Name = "382735__schots__gun-shot";                                Url = "https://freesound.org/people/schots/sounds/382735/" },
    @{ Row = 17; Name = "588246__rkkaleikau__energy-weapon-laser";                 Url = "https://freesound.org/people/rkkaleikau/sounds/588246/" },
    @{ Row = 18; Name = "566435__merrick079__punch2";                             Url = "https://freesound.org/people/Merrick079/sounds/566435/" },
    @{ Row = 19; Name = "232358__richerlandtv__heavy-impacts";                    Url = "https://freesound.org/people/RICHERlandTV/sounds/232358/" },
    @{ Row = 20; Name = "341247__sharesynth__jump01";                             Url = "https://freesound.org/people/sharesynth/sounds/341247/" },
    @{ Row = 21; Name = "561646__mattruthsound__hit-punch-cloth-pillow-bedding-004"; Url = "https://freesound.org/people/MattRuthSound/sounds/561646/" },
    @{ Row = 22; Name = "433644__dersuperanton__game-over-sound";                 Url = "https://freesound.org/people/dersuperanton/sounds/433644/" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("A$row").Value = $r.Name
    $ws.Range("B$row").Value = $r.Url
    $ws.Range("C$row").Value = "Creative Commons License 0"
}

# ---------------------------------------------------------------------
# 3. Hyperlinks: row-insert does not re-anchor the existing Hyperlinks
#    collection, so rebuild it from scratch in the same order Excel wrote
#    it in the target file.
# ---------------------------------------------------------------------
$ws.Hyperlinks.Delete()

$links = @(
    @{ Cell = "B6";  Url = "https://opengameart.org/content/spikes-0" },
    @{ Cell = "B4";  Url = "https://www.fontspace.com/a-area-kilometer-50-font-f53888" },
    @{ Cell = "B2";  Url = "https://freesound.org/people/Whiprealgood/sounds/87535/" },
    @{ Cell = "B3";  Url = "https://freesound.org/people/suntemple/sounds/253172/" },
    @{ Cell = "B5";  Url = "https://opengameart.org/content/simple-explosion-bleeds-game-art" },
    @{ Cell = "B7";  Url = "https://opengameart.org/content/various-inventory-24-pixel-icon-set" },
    @{ Cell = "B8";  Url = "https://opengameart.org/content/energy-icon" },
    @{ Cell = "B26"; Url = "https://free-game-assets.itch.io/night-city-street-2d-background-tiles" },
    @{ Cell = "B25"; Url = "https://elthen.itch.io/2d-pixel-art-vegetable-monsters-sprite-pack" },
    @{ Cell = "B15"; Url = "https://freesound.org/people/qubodup/sounds/60013/" },
    @{ Cell = "B16"; Url = "https://freesound.org/people/schots/sounds/382735/" },
    @{ Cell = "B17"; Url = "https://freesound.org/people/rkkaleikau/sounds/588246/" },
    @{ Cell = "B18"; Url = "https://freesound.org/people/Merrick079/sounds/566435/" },
    @{ Cell = "B21"; Url = "https://freesound.org/people/MattRuthSound/sounds/561646/" },
    @{ Cell = "B20"; Url = "https://freesound.org/people/sharesynth/sounds/341247/" },
    @{ Cell = "B19"; Url = "https://freesound.org/people/RICHERlandTV/sounds/232358/" },
    @{ Cell = "B22"; Url = "https://freesound.org/people/dersuperanton/sounds/433644/" }
)

foreach ($l in $links) {
    $ws.Hyperlinks.Add($ws.Range($l.Cell), $l.Url) | Out-Null
    $ws.Range($l.Cell).Style = "Hyperlink"
}

# ---------------------------------------------------------------------
# 4. View state: active cell / top-left cell after the edit.
# ---------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 16
$ws.Range("C30").Select()
